$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.929813666666668
$ws.Range("H2").Value = 26.789441
$ws.Range("I2").Value = 0.3579859341865942
$ws.Range("J2").Value = 0.3579859341865942
$ws.Range("M2").Value = 50.86142466666666
$ws.Range("N2").Value = 152.584274
$ws.Range("O2").Value = 0.3434314568613803
$ws.Range("P2").Value = 0.3434314568613804
$ws.Range("Q2").Value = 454.1830450945371
$ws.Range("R2").Value = 4087.647405850834
$ws.Range("S2").Value = 0.1229436309135843
$ws.Range("T2").Value = 0.1229436309135843

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.929813666666668
$ws.Range("H3").Value = 26.789441
$ws.Range("I3").Value = 0.3579859341865942
$ws.Range("J3").Value = 0.3579859341865942
$ws.Range("M3").Value = 43.683024
$ws.Range("N3").Value = 131.049072
$ws.Range("O3").Value = 0.294960761928139
$ws.Range("P3").Value = 0.294960761928139
$ws.Range("Q3").Value = 390.081264716528
$ws.Range("R3").Value = 3510.731382448752
$ws.Range("S3").Value = 0.1055918039072344
$ws.Range("T3").Value = 0.1055918039072345

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.929813666666668
$ws.Range("H4").Value = 26.789441
$ws.Range("I4").Value = 0.3579859341865942
$ws.Range("J4").Value = 0.3579859341865942
$ws.Range("M4").Value = 36.64360566666667
$ws.Range("N4").Value = 109.930817
$ws.Range("O4").Value = 0.2474285170192034
$ws.Range("P4").Value = 0.2474285170192035
$ws.Range("Q4").Value = 327.2205706781442
$ws.Range("R4").Value = 2944.985136103297
$ws.Range("S4").Value = 0.08857592880952316
$ws.Range("T4").Value = 0.08857592880952317

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.929813666666668
$ws.Range("H5").Value = 26.789441
$ws.Range("I5").Value = 0.3579859341865942
$ws.Range("J5").Value = 0.3579859341865942
$ws.Range("M5").Value = 16.90969166666667
$ws.Range("N5").Value = 50.729075
$ws.Range("O5").Value = 0.1141792641912772
$ws.Range("P5").Value = 0.1141792641912772
$ws.Range("Q5").Value = 151.0003957441195
$ws.Range("R5").Value = 1359.003561697075
$ws.Range("S5").Value = 0.04087457055625231
$ws.Range("T5").Value = 0.04087457055625231

# Row 6
$ws.Range("I6").Value = 0.03080543264277933
$ws.Range("J6").Value = 0.03080543264277933
$ws.Range("M6").Value = 50.86142466666666
$ws.Range("N6").Value = 152.584274
$ws.Range("O6").Value = 0.3434314568613803
$ws.Range("P6").Value = 0.3434314568613804
$ws.Range("Q6").Value = 39.083393695182
$ws.Range("R6").Value = 351.750543256638
$ws.Range("S6").Value = 0.01057955461175483
$ws.Range("T6").Value = 0.01057955461175483

# Row 7
$ws.Range("I7").Value = 0.03080543264277933
$ws.Range("J7").Value = 0.03080543264277933
$ws.Range("M7").Value = 43.683024
$ws.Range("N7").Value = 131.049072
$ws.Range("O7").Value = 0.294960761928139
$ws.Range("P7").Value = 0.294960761928139
$ws.Range("Q7").Value = 33.56730244929599
$ws.Range("R7").Value = 302.1057220436639
$ws.Range("S7").Value = 0.009086393883840155
$ws.Range("T7").Value = 0.009086393883840157

# Row 8
$ws.Range("I8").Value = 0.03080543264277933
$ws.Range("J8").Value = 0.03080543264277933
$ws.Range("M8").Value = 36.64360566666667
$ws.Range("N8").Value = 109.930817
$ws.Range("O8").Value = 0.2474285170192034
$ws.Range("P8").Value = 0.2474285170192035
$ws.Range("Q8").Value = 28.158009258831
$ws.Range("R8").Value = 253.422083329479
$ws.Range("S8").Value = 0.007622142514937851
$ws.Range("T8").Value = 0.007622142514937852

# Row 9
$ws.Range("I9").Value = 0.03080543264277933
$ws.Range("J9").Value = 0.03080543264277933
$ws.Range("M9").Value = 16.90969166666667
$ws.Range("N9").Value = 50.729075
$ws.Range("O9").Value = 0.1141792641912772
$ws.Range("P9").Value = 0.1141792641912772
$ws.Range("Q9").Value = 12.993897457725
$ws.Range("R9").Value = 116.945077119525
$ws.Range("S9").Value = 0.003517341632246496
$ws.Range("T9").Value = 0.003517341632246497

# Row 10
$ws.Range("G10").Value = 15.246351
$ws.Range("H10").Value = 45.739053
$ws.Range("I10").Value = 0.6112086331706265
$ws.Range("J10").Value = 0.6112086331706265
$ws.Range("M10").Value = 50.86142466666666
$ws.Range("N10").Value = 152.584274
$ws.Range("O10").Value = 0.3434314568613803
$ws.Range("P10").Value = 0.3434314568613804
$ws.Range("Q10").Value = 775.4511328280579
$ws.Range("R10").Value = 6979.060195452522
$ws.Range("S10").Value = 0.2099082713360412
$ws.Range("T10").Value = 0.2099082713360413

# Row 11
$ws.Range("G11").Value = 15.246351
$ws.Range("H11").Value = 45.739053
$ws.Range("I11").Value = 0.6112086331706265
$ws.Range("J11").Value = 0.6112086331706265
$ws.Range("M11").Value = 43.683024
$ws.Range("N11").Value = 131.049072
$ws.Range("O11").Value = 0.294960761928139
$ws.Range("P11").Value = 0.294960761928139
$ws.Range("Q11").Value = 666.0067166454239
$ws.Range("R11").Value = 5994.060449808816
$ws.Range("S11").Value = 0.1802825641370644
$ws.Range("T11").Value = 0.1802825641370644

# Row 12
$ws.Range("G12").Value = 15.246351
$ws.Range("H12").Value = 45.739053
$ws.Range("I12").Value = 0.6112086331706265
$ws.Range("J12").Value = 0.6112086331706265
$ws.Range("M12").Value = 36.64360566666667
$ws.Range("N12").Value = 109.930817
$ws.Range("O12").Value = 0.2474285170192034
$ws.Range("P12").Value = 0.2474285170192035
$ws.Range("Q12").Value = 558.681273899589
$ws.Range("R12").Value = 5028.131465096301
$ws.Range("S12").Value = 0.1512304456947424
$ws.Range("T12").Value = 0.1512304456947424

# Row 13
$ws.Range("G13").Value = 15.246351
$ws.Range("H13").Value = 45.739053
$ws.Range("I13").Value = 0.6112086331706265
$ws.Range("J13").Value = 0.6112086331706265
$ws.Range("M13").Value = 16.90969166666667
$ws.Range("N13").Value = 50.729075
$ws.Range("O13").Value = 0.1141792641912772
$ws.Range("P13").Value = 0.1141792641912772
$ws.Range("Q13").Value = 257.811094451775
$ws.Range("R13").Value = 2320.299850065975
$ws.Range("S13").Value = 0.0697873520027784
$ws.Range("T13").Value = 0.0697873520027784
